$wb = $excel.ActiveWorkbook

# Rename the sheet tabs: drop the "sine_" prefix and rework the ro_ccm/ro_flash
# tokens to ro_CCM / ro_FLASH.
$wb.Worksheets.Item(1).Name = "ro_CCM code_FLASH"
$wb.Worksheets.Item(2).Name = "ro_CCM code_CCM"
$wb.Worksheets.Item(3).Name = "ro_FLASH code_FLASH"
$wb.Worksheets.Item(4).Name = "ro_FLASH code_CCM"

# Update the measured "intensity" (row 2) and "energy" (row 5) values on each
# sheet to the new values recorded for the compare tab.
$values = @{
    1 = @{
        "B2" = 14859;  "C2" = 26277;  "D2" = 33086
        "B5" = 16.535; "C5" = 15.758; "D5" = 17.004
    }
    2 = @{
        "B2" = 12306;  "C2" = 23999;  "D2" = 35495
        "B5" = 13.689; "C5" = 13.365; "D5" = 13.196
    }
    3 = @{
        "B2" = 15337;  "C2" = 24212;  "D2" = 31007
        "B5" = 17.069; "C5" = 18.641; "D5" = 20.366
    }
    4 = @{
        "B2" = 12964;  "C2" = 24385;  "D2" = 34820
        "B5" = 14.42;  "C5" = 14.62;  "D5" = 14.933
    }
}

foreach ($sheetIndex in $values.Keys) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    $cellValues = $values[$sheetIndex]
    foreach ($addr in $cellValues.Keys) {
        $ws.Range($addr).Value = $cellValues[$addr]
    }
}
